$wb = $excel.ActiveWorkbook

# --- Hoja1: update the conversion summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.46 = 50267.85 pesos`n✅ 50267.85 pesos = 12.41 = 973.27 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- tasas: update the N10/O10 and N12/O12 rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 80.27
$ws2.Range("O10").Value = 4035

$ws2.Range("N12").Value = 4051.8
$ws2.Range("O12").Value = 78.45
